# DeveloperGuide: Add retag sequence diagram, update untag and retrieve
# sequence diagrams.
#
# This script only touches the existing "untag sequence diagram" slide
# (slide 1): it relabels the "u:Untag Command" lifeline header to split
# "u:" from "Untag" and widens the two small activation labels that read
# "u" so they read "ut" instead (matching the updated lifeline alias).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1. "Rectangle 5" (the u:Untag Command lifeline box) ---------------
# First paragraph goes from a single run "u:Untag" to two runs:
# "ut:" + "Untag" (the "Command" paragraph below is untouched).
$lifeline = $s.Shapes.Item(3)
$lifelineFirstRun = $lifeline.TextFrame.TextRange.Paragraphs(1).Runs(1)
$lifelineFirstRun.Text = "Untag"
[void]$lifelineFirstRun.InsertBefore("ut:")

# --- 2. "TextBox 44" (activation label "u" near 3135419,2964953) -------
$act1 = $s.Shapes.Item(18)
$act1.TextFrame.TextRange.Text = "ut"
$act1.Width = 421919 / 914400 * 72
$act1.Height = 369332 / 914400 * 72

# --- 3. "TextBox 68" (activation label "u" near 4915208,2617236) -------
$act2 = $s.Shapes.Item(31)
$act2.TextFrame.TextRange.Text = "ut"
$act2.Width = 435921 / 914400 * 72
$act2.Height = 369332 / 914400 * 72
